$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Report Generated On" timestamp
$ws.Range("D5").Value = "Report Generated On: 08/18/2025 09:49 PM"

# Zero out the billed amount / pricing totals (no-violation / no-billable scenario)
$ws.Range("C8").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("H17").Value = 0
